$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("papers")

# Update the pending "Care-seeking behaviour..." row (currently row 163) with its real publication data
$ws.Range("F163").Value = 20
$ws.Range("H163").Value = 260
$ws.Range("K163").Value = "2021-06-09"
$ws.Range("N163").Value = "10.1186/s12936-021-03789-w"
$ws.Range("U163").Value = "OK"

# Promote its formatting to the "processed" look (copy format from the row above)
$ws.Range("A155:U155").Copy()
$ws.Range("A163:U163").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-sort the table by publication date (column K) ascending
$rng = $ws.Range("A1:U163")
$key = $ws.Range("K1")
$rng.Sort($key, 1, $null, $null, 1, $null, 1, 1)

Write-Host "done"
